$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1314.7164
$ws.Range("I15").Value = 1314.7164
$ws.Range("K15").Value = 3944.1492
$ws.Range("M15").Value = -3775.1492
$ws.Range("H26").Value = 2000
$ws.Range("J26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("N26").Value = -2688
$ws.Range("H68").Value = 114420
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 114420
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 114420
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -115918
$ws.Range("H71").Value = 114420
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 114420
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 343260
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -350748
$ws.Range("H113").Value = 2464.375
$ws.Range("I113").Value = 2388.25
$ws.Range("K113").Value = 2388.25
$ws.Range("M113").Value = 865.75
$ws.Range("H116").Value = 58300
$ws.Range("J116").Value = 12450
$ws.Range("L116").Value = 12450
$ws.Range("N116").Value = -19334
$ws.Range("H132").Value = 3578.4348
$ws.Range("I132").Value = 3615.225
$ws.Range("J132").Value = 3333.1667
$ws.Range("K132").Value = 10845.675
$ws.Range("L132").Value = 9999.500100000001
$ws.Range("M132").Value = -8315.674999999999
$ws.Range("N132").Value = -15059.5001
$ws.Range("H141").Value = 4865.7856
$ws.Range("I141").Value = 4177.5835
$ws.Range("K141").Value = 12532.7505
$ws.Range("M141").Value = -7352.750499999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1464.4166
$ws.Range("I2").Value = 932.62964
$ws.Range("J2").Value = 3059.7778
$ws.Range("K2").Value = 932.62964
$ws.Range("L2").Value = 3059.7778
$ws.Range("M2").Value = -819.62964
$ws.Range("N2").Value = -3285.7778
$ws.Range("H32").Value = 16668248
$ws.Range("I32").Value = 6946091.5
$ws.Range("K32").Value = 6946091.5
$ws.Range("M32").Value = -6945804.5
$ws.Range("H52").Value = 74768.5
$ws.Range("J52").Value = 74768.5
$ws.Range("L52").Value = 74768.5
$ws.Range("N52").Value = -75404.5
$ws.Range("H97").Value = 1113.7307
$ws.Range("I97").Value = 772.1818
$ws.Range("K97").Value = 772.1818
$ws.Range("M97").Value = -276.1818
$ws.Range("H116").Value = 1464.4166
$ws.Range("I116").Value = 932.62964
$ws.Range("J116").Value = 3059.7778
$ws.Range("K116").Value = 932.62964
$ws.Range("L116").Value = 3059.7778
$ws.Range("M116").Value = 1361.37036
$ws.Range("N116").Value = -7647.7778
$ws.Range("H122").Value = 2532.7026
$ws.Range("I122").Value = 2109.2058
$ws.Range("J122").Value = 7332.3335
$ws.Range("K122").Value = 6327.617400000001
$ws.Range("L122").Value = 21997.0005
$ws.Range("M122").Value = -3877.617400000001
$ws.Range("N122").Value = -26897.0005
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H132").Value = 2820.2917
$ws.Range("I132").Value = 2454.3057
$ws.Range("J132").Value = 3918.25
$ws.Range("K132").Value = 7362.9171
$ws.Range("L132").Value = 11754.75
$ws.Range("M132").Value = -4832.9171
$ws.Range("N132").Value = -16814.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1464.4166
$ws.Range("I3").Value = 932.62964
$ws.Range("J3").Value = 3059.7778
$ws.Range("K3").Value = 932.62964
$ws.Range("L3").Value = 3059.7778
$ws.Range("M3").Value = -818.62964
$ws.Range("N3").Value = -3287.7778
$ws.Range("H134").Value = 14994159
$ws.Range("I134").Value = 3761816.8
$ws.Range("K134").Value = 11285450.4
$ws.Range("M134").Value = -11282915.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 142857940
$ws.Range("I12").Value = 826
$ws.Range("J12").Value = 333334100
$ws.Range("K12").Value = 826
$ws.Range("L12").Value = 333334100
$ws.Range("M12").Value = -656
$ws.Range("N12").Value = -333334440
$ws.Range("H31").Value = 4350.074
$ws.Range("I31").Value = 2548
$ws.Range("K31").Value = 2548
$ws.Range("M31").Value = -2253
$ws.Range("H34").Value = 4350.074
$ws.Range("I34").Value = 2548
$ws.Range("K34").Value = 2548
$ws.Range("M34").Value = -2346
$ws.Range("H57").Value = 39450
$ws.Range("J57").Value = 39450
$ws.Range("L57").Value = 39450
$ws.Range("N57").Value = -40570
$ws.Range("H122").Value = 2950.7646
$ws.Range("I122").Value = 2436
$ws.Range("J122").Value = 4186.2
$ws.Range("K122").Value = 7308
$ws.Range("L122").Value = 12558.6
$ws.Range("M122").Value = -4858
$ws.Range("N122").Value = -17458.6
$ws.Range("H127").Value = 111938.6
$ws.Range("J127").Value = 114923
$ws.Range("L127").Value = 114923
$ws.Range("N127").Value = -124843
$ws.Range("H134").Value = 2770.2693
$ws.Range("I134").Value = 1422.6842
$ws.Range("K134").Value = 4268.0526
$ws.Range("M134").Value = -1733.0526

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62916484
$ws.Range("I4").Value = 37941904
$ws.Range("K4").Value = 113825712
$ws.Range("M4").Value = -113825600
$ws.Range("H34").Value = 584
$ws.Range("J34").Value = 1500
$ws.Range("L34").Value = 4500
$ws.Range("N34").Value = -4668
$ws.Range("H39").Value = 2333.3333
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 6000
$ws.Range("M39").Value = -5706
$ws.Range("H55").Value = 2526
$ws.Range("I55").Value = 1902
$ws.Range("J55").Value = 3150
$ws.Range("K55").Value = 5706
$ws.Range("L55").Value = 9450
$ws.Range("M55").Value = -5529
$ws.Range("N55").Value = -9804
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2730
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2064
$ws.Range("H109").Value = 1996.6666
$ws.Range("I109").Value = 1996.6666
$ws.Range("K109").Value = 5989.9998
$ws.Range("M109").Value = -4949.9998
$ws.Range("H132").Value = 1462.64
$ws.Range("I132").Value = 934.1111
$ws.Range("J132").Value = 1759.9375
$ws.Range("K132").Value = 8406.999899999999
$ws.Range("L132").Value = 15839.4375
$ws.Range("M132").Value = -5876.999899999999
$ws.Range("N132").Value = -20899.4375
$ws.Range("H140").Value = 20002570
$ws.Range("I140").Value = 20002570
$ws.Range("K140").Value = 60007710
$ws.Range("M140").Value = -60002530

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3727.476
$ws.Range("I132").Value = 3763.85
$ws.Range("K132").Value = 11291.55
$ws.Range("M132").Value = -8761.549999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1931
$ws.Range("I61").Value = 1931
$ws.Range("K61").Value = 1931
$ws.Range("M61").Value = -1729
$ws.Range("H68").Value = 2388.7144
$ws.Range("I68").Value = 1180.25
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1180.25
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -431.25
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2388.7144
$ws.Range("I71").Value = 1180.25
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 5901.25
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -2157.25
$ws.Range("N71").Value = -27488
$ws.Range("H113").Value = 1931
$ws.Range("I113").Value = 1931
$ws.Range("K113").Value = 1931
$ws.Range("M113").Value = 239
$ws.Range("H122").Value = 3031.75
$ws.Range("I122").Value = 2991.0454
$ws.Range("K122").Value = 8973.136200000001
$ws.Range("M122").Value = -6523.136200000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 44350
$ws.Range("I96").Value = 27500
$ws.Range("K96").Value = 27500
$ws.Range("M96").Value = -26127
$ws.Range("H122").Value = 2194.8965
$ws.Range("I122").Value = 1793
$ws.Range("J122").Value = 4706.75
$ws.Range("K122").Value = 5379
$ws.Range("L122").Value = 14120.25
$ws.Range("M122").Value = -2929
$ws.Range("N122").Value = -19020.25
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H130").Value = 88876
$ws.Range("J130").Value = 88876
$ws.Range("L130").Value = 88876
$ws.Range("N130").Value = -98916
$ws.Range("H132").Value = 3489.913
$ws.Range("I132").Value = 2815.3157
$ws.Range("K132").Value = 8445.947100000001
$ws.Range("M132").Value = -5915.947100000001
